$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.687.25'
$ws.Range("E2").Value = '  +0.34%  '
$ws.Range("D3").Value = '2.459.12'
$ws.Range("E3").Value = '  +0.28%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '''559.39'
$ws.Range("E5").Value = '  -0.88%  '
$ws.Range("D6").Value = '''161.66'
$ws.Range("E6").Value = '  -1.36%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '''0.505'
$ws.Range("E8").Value = '  -0.25%  '
$ws.Range("E9").Value = '  -0.50%  '
$ws.Range("E10").Value = '  +0.52%  '
$ws.Range("E11").Value = '  -2.88%  '
$ws.Range("E12").Value = '  +0.79%  '
$ws.Range("E13").Value = '  +0.11%  '
$ws.Range("D14").Value = '68.591.05'
$ws.Range("E14").Value = '  +0.30%  '
$ws.Range("E15").Value = '  -1.84%  '
$ws.Range("D16").Value = '''23.51'
$ws.Range("E16").Value = '  -0.44%  '
$ws.Range("D17").Value = '2.437.59'
$ws.Range("E17").Value = '  -0.94%  '
$ws.Range("D18").Value = '''10.64'
$ws.Range("E18").Value = '  -3.23%  '
$ws.Range("D19").Value = '''334.09'
$ws.Range("E19").Value = '  -3.11%  '
$ws.Range("E20").Value = '  -3.59%  '
$ws.Range("E21").Value = '  -1.37%  '
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("E23").Value = '  -0.66%  '
$ws.Range("D24").Value = '''66.45'
$ws.Range("E24").Value = '  -2.50%  '
$ws.Range("E25").Value = '  -3.42%  '
$ws.Range("E26").Value = '  -1.54%  '
$ws.Range("E27").Value = '  -3.63%  '
$ws.Range("E28").Value = '  -2.10%  '
$ws.Range("E29").Value = '  -0.05%  '
$ws.Range("D30").Value = '''429.01'
$ws.Range("E30").Value = '  -1.96%  '
$ws.Range("E31").Value = '  -4.35%  '
$ws.Range("E32").Value = '  -4.63%  '
$ws.Range("D33").Value = '''158.91'
$ws.Range("E33").Value = '  +1.34%  '
$ws.Range("D34").Value = '''19.01'
$ws.Range("E34").Value = '  +0.08%  '
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("E36").Value = '  -0.61%  '
$ws.Range("E37").Value = '  -1.11%  '
$ws.Range("E38").Value = '  -2.51%  '
$ws.Range("D39").Value = '''4.40'
$ws.Range("E39").Value = '  -2.20%  '
$ws.Range("E40").Value = '  -4.88%  '
$ws.Range("E41").Value = '  -6.31%  '
$ws.Range("D42").Value = '''2.06'
$ws.Range("E42").Value = '  -1.70%  '
$ws.Range("E43").Value = '  -0.96%  '
$ws.Range("D44").Value = '''129.89'
$ws.Range("E44").Value = '  -4.01%  '
$ws.Range("E45").Value = '  -0.68%  '
$ws.Range("D46").Value = '''0.481'
$ws.Range("E46").Value = '  -1.33%  '
$ws.Range("D47").Value = '''0.558'
$ws.Range("E47").Value = '  -1.25%  '
$ws.Range("E48").Value = '  -0.98%  '
$ws.Range("E49").Value = '  +0.17%  '
$ws.Range("E50").Value = '  -3.66%  '
$ws.Range("D51").Value = '''4.95'
$ws.Range("E51").Value = '  -8.55%  '
